$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 50, shifting rows 50-53 down to 51-54
$ws.Rows.Item(50).Insert()

# New row 50 values (Dina / Primera / 44931)
$ws.Range("A50").Value = 5
$ws.Range("B50").Value = "Macroferia Regional de Talca"
$ws.Range("C50").Value = "Maule"
$ws.Range("D50").Value = 44931
$ws.Range("E50").Value = 7
$ws.Range("F50").Value = "Fruta"
$ws.Range("G50").Value = 100103
$ws.Range("H50").Value = "Frutos de hueso (carozo)"
$ws.Range("I50").Value = 100103003
$ws.Range("J50").Value = "Damasco"
$ws.Range("K50").Value = "Dina"
$ws.Range("L50").Value = "Primera"
$ws.Range("M50").Value = 180
$ws.Range("N50").Value = 15000
$ws.Range("O50").Value = 15000
$ws.Range("P50").Value = 15000
$ws.Range("Q50").Value = "`$/caja 16 kilos"
$ws.Range("R50").Value = "Región de O'Higgins"
$ws.Range("S50").Value = 938
$ws.Range("T50").Value = 16
